$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2222222222222222
$ws.Range("C2").Value = 0.4444444444444444
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("P3").Value = 1
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("D6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.1818181818181818
$ws.Range("Q6").Value = 0.3636363636363636
$ws.Range("R6").Value = 0.1818181818181818
$ws.Range("S6").Value = 0.1818181818181818
$ws.Range("Q7").Value = 0.3333333333333333
$ws.Range("R7").Value = 0.1666666666666667
$ws.Range("S7").Value = 0.5
$ws.Range("D8").Value = 0.03846153846153846
$ws.Range("J8").Value = 0.1923076923076923
$ws.Range("O8").Value = 0.03846153846153846
$ws.Range("Q8").Value = 0.1153846153846154
$ws.Range("R8").Value = 0.1538461538461539
$ws.Range("S8").Value = 0.4615384615384616
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.8571428571428571
$ws.Range("B10").Value = 0.05084745762711865
$ws.Range("D10").Value = 0.05084745762711865
$ws.Range("E10").Value = 0.008474576271186441
$ws.Range("F10").Value = 0.02542372881355932
$ws.Range("J10").Value = 0.1271186440677966
$ws.Range("O10").Value = 0.0423728813559322
$ws.Range("Q10").Value = 0.288135593220339
$ws.Range("R10").Value = 0.06779661016949153
$ws.Range("S10").Value = 0.3389830508474576
$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.1111111111111111
$ws.Range("L11").Value = 0.7777777777777778
$ws.Range("G12").Value = 0.5555555555555556
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("L12").Value = 0.2222222222222222
$ws.Range("F15").Value = 0.04166666666666666
$ws.Range("H15").Value = 0.08333333333333333
$ws.Range("I15").Value = 0.04166666666666666
$ws.Range("O15").Value = 0.125
$ws.Range("S15").Value = 0.2083333333333333
$ws.Range("F16").Value = 0.07692307692307693
$ws.Range("H16").Value = 0.1538461538461539
$ws.Range("J16").Value = 0.7692307692307693
$ws.Range("F17").Value = 0.04761904761904762
$ws.Range("H17").Value = 0.119047619047619
$ws.Range("I17").Value = 0.04761904761904762
$ws.Range("J17").Value = 0.6190476190476191
$ws.Range("K17").Value = 0.02380952380952381
$ws.Range("M17").Value = 0.02380952380952381
$ws.Range("O17").Value = 0.04761904761904762
$ws.Range("S17").Value = 0.07142857142857142
$ws.Range("H18").Value = 0.0625
$ws.Range("J18").Value = 0.625
$ws.Range("K18").Value = 0.125
$ws.Range("O18").Value = 0.125
$ws.Range("S18").Value = 0.0625
$ws.Range("F19").Value = 0.04
$ws.Range("H19").Value = 0.2133333333333333
$ws.Range("I19").Value = 0.04
$ws.Range("J19").Value = 0.5066666666666667
$ws.Range("K19").Value = 0.05333333333333334
$ws.Range("O19").Value = 0.09333333333333334
$ws.Range("S19").Value = 0.05333333333333334
